$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 223.83333
$ws.Range("I38").Value = 223.83333
$ws.Range("K38").Value = 671.49999
$ws.Range("M38").Value = -299.49999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H49").Value = 1000
$ws.Range("J49").Value = 1000
$ws.Range("L49").Value = 3000
$ws.Range("N49").Value = -3272

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 73089.45
$ws.Range("I62").Value = 79398.5
$ws.Range("J62").Value = 9999
$ws.Range("K62").Value = 79398.5
$ws.Range("L62").Value = 9999
$ws.Range("M62").Value = -78774.5
$ws.Range("N62").Value = -11247

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 73089.45
$ws.Range("I65").Value = 79398.5
$ws.Range("J65").Value = 9999
$ws.Range("K65").Value = 396992.5
$ws.Range("L65").Value = 49995
$ws.Range("M65").Value = -393872.5
$ws.Range("N65").Value = -56235

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 6924.2246
$ws.Range("I132").Value = 5132.619
$ws.Range("K132").Value = 15397.857
$ws.Range("M132").Value = -12867.857

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3039.7036
$ws.Range("I137").Value = 2443.9333
$ws.Range("J137").Value = 3784.4167
$ws.Range("K137").Value = 7331.7999
$ws.Range("L137").Value = 11353.2501
$ws.Range("M137").Value = -4781.7999
$ws.Range("N137").Value = -16453.2501

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 7359.6562
$ws.Range("I138").Value = 9565.916999999999
$ws.Range("J138").Value = 6850.519
$ws.Range("K138").Value = 28697.751
$ws.Range("L138").Value = 20551.557
$ws.Range("M138").Value = -23557.751
$ws.Range("N138").Value = -30831.557

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1886
$ws.Range("I2").Value = 1941.1428
$ws.Range("K2").Value = 1941.1428
$ws.Range("M2").Value = -1828.1428

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 41966.35
$ws.Range("I32").Value = 24776.857
$ws.Range("J32").Value = 53999
$ws.Range("K32").Value = 24776.857
$ws.Range("L32").Value = 53999
$ws.Range("M32").Value = -24489.857
$ws.Range("N32").Value = -54573

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6754.125
$ws.Range("I61").Value = 6603.615
$ws.Range("J61").Value = 7406.3335
$ws.Range("K61").Value = 6603.615
$ws.Range("L61").Value = 7406.3335
$ws.Range("M61").Value = -6391.615
$ws.Range("N61").Value = -7830.3335

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H103").Value = 95000
$ws.Range("J103").Value = 95000
$ws.Range("L103").Value = 95000
$ws.Range("N103").Value = -97344

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1886
$ws.Range("I116").Value = 1941.1428
$ws.Range("K116").Value = 1941.1428
$ws.Range("M116").Value = 352.8571999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 27617.586
$ws.Range("I132").Value = 32478.646
$ws.Range("J132").Value = 4006.7144
$ws.Range("K132").Value = 97435.93799999999
$ws.Range("L132").Value = 12020.1432
$ws.Range("M132").Value = -94905.93799999999
$ws.Range("N132").Value = -17080.1432

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 6754.125
$ws.Range("I136").Value = 6603.615
$ws.Range("J136").Value = 7406.3335
$ws.Range("K136").Value = 19810.845
$ws.Range("L136").Value = 22219.0005
$ws.Range("M136").Value = -17260.845
$ws.Range("N136").Value = -27319.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1886
$ws.Range("I3").Value = 1941.1428
$ws.Range("K3").Value = 1941.1428
$ws.Range("M3").Value = -1827.1428

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3024
$ws.Range("I31").Value = 2630.3333
$ws.Range("J31").Value = 4598.6665
$ws.Range("K31").Value = 2630.3333
$ws.Range("L31").Value = 4598.6665
$ws.Range("M31").Value = -2335.3333
$ws.Range("N31").Value = -5188.6665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3024
$ws.Range("I34").Value = 2630.3333
$ws.Range("J34").Value = 4598.6665
$ws.Range("K34").Value = 2630.3333
$ws.Range("L34").Value = 4598.6665
$ws.Range("M34").Value = -2428.3333
$ws.Range("N34").Value = -5002.6665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H44").Value = 6500
$ws.Range("I44").Value = 6500
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 6500
$ws.Range("L44").Value = 0
$ws.Range("N44").Value = -6058
$ws.Range("M44").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 28861
$ws.Range("I134").Value = 33043.47
$ws.Range("J134").Value = 5160.3335
$ws.Range("K134").Value = 99130.41
$ws.Range("L134").Value = 15481.0005
$ws.Range("M134").Value = -96595.41
$ws.Range("N134").Value = -20551.0005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 453543.9
$ws.Range("J141").Value = 469871.9
$ws.Range("L141").Value = 469871.9
$ws.Range("N141").Value = -480231.9

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 5888.222
$ws.Range("J80").Value = 5888.222
$ws.Range("L80").Value = 17664.666
$ws.Range("N80").Value = -19536.666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value = 5888.222
$ws.Range("J83").Value = 5888.222
$ws.Range("L83").Value = 52993.998
$ws.Range("N83").Value = -62353.998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 460.16666
$ws.Range("J107").Value = 682.3333
$ws.Range("L107").Value = 2046.9999
$ws.Range("N107").Value = -5886.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 3131.6
$ws.Range("I113").Value = 1500
$ws.Range("J113").Value = 3539.5
$ws.Range("K113").Value = 4500
$ws.Range("L113").Value = 10618.5
$ws.Range("M113").Value = -2330
$ws.Range("N113").Value = -14958.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3797.4
$ws.Range("I80").Value = 3495.6667
$ws.Range("J80").Value = 4250
$ws.Range("K80").Value = 3495.6667
$ws.Range("L80").Value = 4250
$ws.Range("M80").Value = -2497.6667
$ws.Range("N80").Value = -6246

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3797.4
$ws.Range("I83").Value = 3495.6667
$ws.Range("J83").Value = 4250
$ws.Range("K83").Value = 17478.3335
$ws.Range("L83").Value = 21250
$ws.Range("M83").Value = -12486.3335
$ws.Range("N83").Value = -31234

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2000
$ws.Range("I97").Value = 2000
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 2000
$ws.Range("L97").Value = 0
$ws.Range("N97").Value = -1504
$ws.Range("M97").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 58389.39
$ws.Range("I113").Value = 74273.14
$ws.Range("J113").Value = 2796.25
$ws.Range("K113").Value = 74273.14
$ws.Range("L113").Value = 2796.25
$ws.Range("M113").Value = -72103.14
$ws.Range("N113").Value = -7136.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3479.6365
$ws.Range("I122").Value = 3327.6
$ws.Range("K122").Value = 9982.799999999999
$ws.Range("M122").Value = -7532.799999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 738.03125
$ws.Range("I55").Value = 489.21054
$ws.Range("K55").Value = 489.21054
$ws.Range("M55").Value = -316.21054

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H97").Value = 19074.75
$ws.Range("J97").Value = 19074.75
$ws.Range("L97").Value = 19074.75
$ws.Range("N97").Value = -21056.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 36143.566
$ws.Range("I132").Value = 43487.832
$ws.Range("K132").Value = 130463.496
$ws.Range("M132").Value = -127933.496

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H117").Value = 20000
$ws.Range("J117").Value = 20000
$ws.Range("L117").Value = 20000
$ws.Range("N117").Value = -29178

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3481.6365
$ws.Range("I122").Value = 3029.8
$ws.Range("J122").Value = 8000
$ws.Range("K122").Value = 9089.400000000001
$ws.Range("L122").Value = 24000
$ws.Range("M122").Value = -6639.400000000001
$ws.Range("N122").Value = -28900

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 48334.93
$ws.Range("I132").Value = 54208.273
$ws.Range("K132").Value = 162624.819
$ws.Range("M132").Value = -160094.819
